$d = $word.ActiveDocument

# The document ends with a single empty paragraph. We turn it into a bold
# "question" paragraph ("To add my name to the README.md") and add a new,
# non-bold paragraph right after it containing the answer text.

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastRange = $lastPara.Range

# Create the new (second) paragraph right after the current last paragraph,
# inheriting the same (non-bold) paragraph/run formatting as the original
# empty paragraph.
$lastRange.InsertParagraphAfter()

# Paragraph that used to be last now becomes the bold heading line.
$headingPara = $d.Paragraphs($lastIndex)
$headingRange = $headingPara.Range
$headingRange.Text = "To add my name to the README.md"
$headingRange.Font.Bold = 1

# The newly inserted paragraph becomes the descriptive answer line.
$bodyPara = $d.Paragraphs($lastIndex + 1)
$bodyRange = $bodyPara.Range
$bodyRange.Text = "Create a fork of the paceuniversity/courses repository. Edit the forked README.md to add my name. Create a pull request with a comment that says to add my info to the bottom of the file. "
